$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E3").Value = 22
$ws.Range("E7").Value = 25
$ws.Range("E12").Value = 23
$ws.Range("E16").Value = 292
$ws.Range("E18").Value = 88
